# Auto-generated Excel COM-interop script
# Applies the "Update countries & provincias Spain" data refresh to the Pais sheet:
# re-ranks India/Canada, Sudafrica/Dinamarca/Republica Dominicana/Serbia,
# Mozambique/Aruba/Monaco/Bahamas and Belice/Nueva Caledonia, and refreshes their
# case counts (columns B-H) to match the latest daily snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: -> Estados Unidos
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 1370436
$ws.Cells.Item(4, 3).Value = 2798
$ws.Cells.Item(4, 4).Value = 256346
$ws.Cells.Item(4, 5).Value = 1033227
$ws.Cells.Item(4, 6).Value = 16514
$ws.Cells.Item(4, 7).Value = 76
$ws.Cells.Item(4, 8).Value = 80863

# Row 15: -> India
$ws.Cells.Item(15, 1).Value = 'India'
$ws.Cells.Item(15, 2).Value = 69149
$ws.Cells.Item(15, 3).Value = 1988
$ws.Cells.Item(15, 4).Value = 21664
$ws.Cells.Item(15, 5).Value = 45236
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 37
$ws.Cells.Item(15, 8).Value = 2249

# Row 16: -> Canada
$ws.Cells.Item(16, 1).Value = 'Canada'
$ws.Cells.Item(16, 2).Value = 68848
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 32096
$ws.Cells.Item(16, 5).Value = 31882
$ws.Cells.Item(16, 6).Value = 502
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 4870

# Row 23: -> Suiza
$ws.Cells.Item(23, 1).Value = 'Suiza'
$ws.Cells.Item(23, 2).Value = 30344
$ws.Cells.Item(23, 3).Value = 39
$ws.Cells.Item(23, 4).Value = 26600
$ws.Cells.Item(23, 5).Value = 1907
$ws.Cells.Item(23, 6).Value = 89
$ws.Cells.Item(23, 7).Value = 4
$ws.Cells.Item(23, 8).Value = 1837

# Row 44: -> Sudafrica
$ws.Cells.Item(44, 1).Value = 'Sudafrica'
$ws.Cells.Item(44, 2).Value = 10652
$ws.Cells.Item(44, 3).Value = 637
$ws.Cells.Item(44, 4).Value = 4357
$ws.Cells.Item(44, 5).Value = 6089
$ws.Cells.Item(44, 6).Value = 77
$ws.Cells.Item(44, 7).Value = 12
$ws.Cells.Item(44, 8).Value = 206

# Row 45: -> Dinamarca
$ws.Cells.Item(45, 1).Value = 'Dinamarca'
$ws.Cells.Item(45, 2).Value = 10513
$ws.Cells.Item(45, 3).Value = 84
$ws.Cells.Item(45, 4).Value = 8328
$ws.Cells.Item(45, 5).Value = 1652
$ws.Cells.Item(45, 6).Value = 43
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = 533

# Row 46: -> Republica Dominicana
$ws.Cells.Item(46, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(46, 2).Value = 10347
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 2763
$ws.Cells.Item(46, 5).Value = 7196
$ws.Cells.Item(46, 6).Value = 134
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 388

# Row 47: -> Serbia
$ws.Cells.Item(47, 1).Value = 'Serbia'
$ws.Cells.Item(47, 2).Value = 10176
$ws.Cells.Item(47, 3).Value = 62
$ws.Cells.Item(47, 4).Value = 3290
$ws.Cells.Item(47, 5).Value = 6668
$ws.Cells.Item(47, 6).Value = 29
$ws.Cells.Item(47, 7).Value = 3
$ws.Cells.Item(47, 8).Value = 218

# Row 56: -> Argentina
$ws.Cells.Item(56, 1).Value = 'Argentina'
$ws.Cells.Item(56, 2).Value = 6034
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 1837
$ws.Cells.Item(56, 5).Value = 3892
$ws.Cells.Item(56, 6).Value = 148
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 305

# Row 163: -> Mozambique
$ws.Cells.Item(163, 1).Value = 'Mozambique'
$ws.Cells.Item(163, 2).Value = 103
$ws.Cells.Item(163, 3).Value = 12
$ws.Cells.Item(163, 4).Value = 34
$ws.Cells.Item(163, 5).Value = 69
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 0

# Row 164: -> Aruba
$ws.Cells.Item(164, 1).Value = 'Aruba'
$ws.Cells.Item(164, 2).Value = 101
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 89
$ws.Cells.Item(164, 5).Value = 9
$ws.Cells.Item(164, 6).Value = 4
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 3

# Row 165: -> Monaco
$ws.Cells.Item(165, 1).Value = 'Monaco'
$ws.Cells.Item(165, 2).Value = 96
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 82
$ws.Cells.Item(165, 5).Value = 10
$ws.Cells.Item(165, 6).Value = 1
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 4

# Row 166: -> Bahamas
$ws.Cells.Item(166, 1).Value = 'Bahamas'
$ws.Cells.Item(166, 2).Value = 92
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 37
$ws.Cells.Item(166, 5).Value = 44
$ws.Cells.Item(166, 6).Value = 1
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 11

# Row 192: -> Belice
$ws.Cells.Item(192, 1).Value = 'Belice'
$ws.Cells.Item(192, 2).Value = 18
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 16
$ws.Cells.Item(192, 5).Value = 0
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 2

# Row 193: -> Nueva Caledonia
$ws.Cells.Item(193, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(193, 2).Value = 18
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 18
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

